# Add a header row (Serial-No, Emp-id, Name) to Sheet1 and leave the
# selection on D10, matching the authored change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Serial-No"
$ws.Range("B1").Value = "Emp-id"
$ws.Range("C1").Value = "Name"

$ws.Range("D10").Select()
